$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data clean up -------------------------------------------------------
# Insert a summary/header row above the existing recipe table: a label in
# A1 describing the sheet, plus a running count (1..4) above each of the
# four recipe blocks (columns A, C, E, G, I hold Recipe ID / Name / etc.
# for recipes 1-4 respectively).
$ws.Range("A1").Value = "No of Reccipes"
$ws.Range("C1").Value = 1
$ws.Range("E1").Value = 2
$ws.Range("G1").Value = 3
$ws.Range("I1").Value = 4

# Re-apply the column width for A:C (kept at 39 characters, just re-saved)
$ws.Columns.Item(1).ColumnWidth = 38.1666666666667
$ws.Columns.Item(2).ColumnWidth = 38.1666666666667
$ws.Columns.Item(3).ColumnWidth = 38.1666666666667

# --- View state ------------------------------------------------------------
# Scroll the window so column G is in view and move the active selection
# over to Z1, matching where the workbook was left when last saved.
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$win.ScrollRow = 1
$ws.Range("Z1").Select()
